# Fix typo in race name: "10. Łańcucka Piętka" -> "10. Łańcucka Piątka"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("A2").Value = "10. Łańcucka Piątka"

# Update the active selection to match the saved view state (B6)
$ws.Range("B6").Select()
